$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
32,32,11,14,25,31,21,21,40,19,
34,6,20,40,39,20,33,25,39,7,
13,31,31,23,25,53,53,4,45,35,
18,2,28,4,35,22,4,28,30,26,
4,41,39,1,28,18,59,38,17,52,
34,10,5,14,8,2,32,33,48,27,
13,45,49,26,22,15,36,19,10,26,
51,28,30,27,28,13,13,25,15,17,
15,29,16,38,3,15,9,50,35,39,
26,18,11,25,40,43,7,27,30,34,
40,50,37,20,4,9,20,13,35,35,
33,8,10,48,35,31,26,26,7,25,
41,26,4,42,28,38,22,17,40,6,
23,40,22,34,18,16,22,41,8,24,
3,22,42,23,30,22,8,16,5,46,
37,46,3,40,14,26,1,52,23,50
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 4).Value = $values[$i]
}

$ws.Range("E1:E1048576").Select()
